# Auto-generated script: update cryptos price/volume table
# Applies text-safe cell value updates (preserving string formatting,
# e.g. "1.00" stays "1.00" instead of being coerced to the number 1).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $cell = $ws.Range($cellRef)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.ClearFormats()
}

Set-TextValue 'D2' '67.972.27'

Set-TextValue 'D3' '3.259.55'
Set-TextValue 'E3' '  -0.86%  '

Set-TextValue 'E4' '  -0.01%  '

Set-TextValue 'D5' '582.72'
Set-TextValue 'E5' '  -0.16%  '

Set-TextValue 'D6' '184.11'
Set-TextValue 'E6' '  -1.28%  '

Set-TextValue 'E7' '  +0.02%  '

Set-TextValue 'E8' '  -0.06%  '

Set-TextValue 'E9' '  -3.59%  '

Set-TextValue 'E10' '  -0.98%  '

Set-TextValue 'E11' '  -3.68%  '

Set-TextValue 'D12' '3.826.07'
Set-TextValue 'E12' '  -0.73%  '

Set-TextValue 'E13' '  +1.43%  '

Set-TextValue 'B14' 'Avalanche'
Set-TextValue 'C14' 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
Set-TextValue 'D14' '27.35'
Set-TextValue 'E14' '  -4.30%  '

Set-TextValue 'B15' 'WrappedBTC'
Set-TextValue 'C15' 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
Set-TextValue 'D15' '67.973.53'
Set-TextValue 'E15' '  -1.32%  '

Set-TextValue 'E16' '  -2.76%  '

Set-TextValue 'D17' '3.256.03'
Set-TextValue 'E17' '  -0.40%  '

Set-TextValue 'E18' '  -2.64%  '

Set-TextValue 'D19' '13.24'
Set-TextValue 'E19' '  -3.07%  '

Set-TextValue 'D20' '416.71'
Set-TextValue 'E20' '  +5.30%  '

Set-TextValue 'D21' '7.52'
Set-TextValue 'E21' '  -3.23%  '

Set-TextValue 'E22' '  -0.09%  '

Set-TextValue 'D23' '71.06'

Set-TextValue 'D24' '0.507'

Set-TextValue 'E25' '  -3.35%  '

Set-TextValue 'E26' '  -1.41%  '

Set-TextValue 'D27' '9.35'
Set-TextValue 'E27' '  -3.91%  '

Set-TextValue 'D28' '1.00'
Set-TextValue 'E28' '  +0.31%  '

Set-TextValue 'E29' '  -2.09%  '

Set-TextValue 'D30' '22.54'
Set-TextValue 'E30' '  -2.84%  '

Set-TextValue 'E31' '  -5.37%  '

Set-TextValue 'D32' '6.84'
Set-TextValue 'E32' '  -5.02%  '

Set-TextValue 'E33' '  -5.20%  '

Set-TextValue 'D34' '162.98'
Set-TextValue 'E34' '  -0.49%  '

Set-TextValue 'E35' '  -5.69%  '

Set-TextValue 'E36' '  -4.93%  '

Set-TextValue 'D37' '26.81'
Set-TextValue 'E37' '  -0.29%  '

Set-TextValue 'D38' '0.794'
Set-TextValue 'E38' '  -4.22%  '

Set-TextValue 'E39' '  -4.13%  '

Set-TextValue 'E40' '  -5.95%  '

Set-TextValue 'D41' '2.632.71'
Set-TextValue 'E41' '  -0.92%  '

Set-TextValue 'B42' 'dogwifhat'
Set-TextValue 'C42' 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
Set-TextValue 'D42' '2.42'
Set-TextValue 'E42' '  -5.82%  '

Set-TextValue 'B43' 'Hedera'
Set-TextValue 'C43' 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextValue 'D43' '0.0672'
Set-TextValue 'E43' '  -2.83%  '

Set-TextValue 'D44' '337.06'
Set-TextValue 'E44' '  -1.58%  '

Set-TextValue 'D45' '24.21'
Set-TextValue 'E45' '  -5.55%  '

Set-TextValue 'E46' '  -3.89%  '

Set-TextValue 'D47' '6.22'
Set-TextValue 'E47' '  -2.32%  '

Set-TextValue 'B48' 'Stellar'
Set-TextValue 'C48' 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
Set-TextValue 'D48' '0.100'
Set-TextValue 'E48' '  -2.64%  '

Set-TextValue 'B49' 'ONDO'
Set-TextValue 'C49' 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
Set-TextValue 'D49' '0.972'
Set-TextValue 'E49' '  -2.89%  '

Set-TextValue 'D50' '0.999'
Set-TextValue 'E50' '  +0.00%  '

Set-TextValue 'E51' '  -5.27%  '
